# Commit: "Add vehicle maintenance costs (#36)"
#
# The "Key to Variables" sheet gets one new acronym entry inserted right
# before the existing "BAADTbVT" row (old row 153), pushing every row
# below it down by one:
#   Top Level Folder = trans
#   Acronym           = AVMC
#   Meaning           = Annual Vehicle Maintenance Cost
#   Importance to Update for New Country = low
#
# The author also had the "Key to Variables" tab active/selected when the
# file was saved (previously the "About" tab was selected), so we switch
# to that sheet as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# Insert a new blank row above the current row 153 - everything currently
# at row 153 and below shifts down to make room.
$ws.Rows.Item(153).Insert()

# Populate the new row's data.
$ws.Cells.Item(153, 1).Value = "trans"
$ws.Cells.Item(153, 2).Value = "AVMC"
$ws.Cells.Item(153, 3).Value = "Annual Vehicle Maintenance Cost"
$ws.Cells.Item(153, 6).Value = "low"

# The freshly-inserted row picked up the "medium" shading from the row
# above it (default Excel "insert copies formatting from above" behavior).
# Column F should instead carry the "low" shading used elsewhere in the
# sheet, so copy that formatting over from another "low" row.
$ws.Cells.Item(167, 6).Copy()
$ws.Cells.Item(153, 6).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Make "Key to Variables" the active sheet/tab, matching the saved view
# state captured in the workbook.
$ws.Activate()
